$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (masthead volume/number and reporting week dates) ---
$ws.Range("A8").Value = "Volume 29   Number  44"
$ws.Range("C9").Value = "Report Covering the Week  10/31/2022  Through  11/6/2022"

# --- Weekly crime statistics grid updates (rows 15-30) ---
$ws.Range("L15").Value = 140
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -66.666666666666
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = -62.5
$ws.Range("I16").Value = 98
$ws.Range("J16").Value = 95
$ws.Range("K16").Value = 3.157894736842
$ws.Range("L16").Value = 27.272727272727
$ws.Range("M16").Value = 28.947368421052
$ws.Range("N16").Value = -81.153846153846
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -40
$ws.Range("F17").Value = 10
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = -33.333333333333
$ws.Range("I17").Value = 90
$ws.Range("J17").Value = 65
$ws.Range("K17").Value = 38.461538461538
$ws.Range("L17").Value = 87.5
$ws.Range("M17").Value = 80
$ws.Range("N17").Value = 9.756097560975
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = -71.428571428571
$ws.Range("L18").Value = -7.086614173228
$ws.Range("M18").Value = 22.916666666666
$ws.Range("N18").Value = -84.793814432989
$ws.Range("D19").Value = 21
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 73
$ws.Range("G19").Value = 99
$ws.Range("H19").Value = -26.262626262626
$ws.Range("I19").Value = 719
$ws.Range("J19").Value = 628
$ws.Range("K19").Value = 14.490445859872
$ws.Range("L19").Value = 74.09200968523
$ws.Range("M19").Value = 28.853046594982
$ws.Range("N19").Value = -54.232972628898
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = -44.444444444444
$ws.Range("J20").Value = 80
$ws.Range("K20").Value = -18.75
$ws.Range("L20").Value = 66.666666666666
$ws.Range("N20").Value = -93.434343434343
$ws.Range("C21").Value = 25
$ws.Range("D21").Value = 30
$ws.Range("E21").Value = -16.666666666666
$ws.Range("F21").Value = 96
$ws.Range("G21").Value = 146
$ws.Range("H21").Value = -34.246575342465
$ws.Range("I21").Value = 1102
$ws.Range("J21").Value = 942
$ws.Range("K21").Value = 16.985138004246
$ws.Range("L21").Value = 55.211267605633
$ws.Range("M21").Value = 34.718826405868
$ws.Range("N21").Value = -72.164688052538
$ws.Range("D22").Value = 2
$ws.Range("D22").NumberFormat = "#,##0"
$ws.Range("E22").Value = -100
$ws.Range("E22").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -66.666666666666
$ws.Range("J22").Value = 21
$ws.Range("K22").Value = 4.761904761904
$ws.Range("L22").Value = 29.411764705882
$ws.Range("M22").Value = 0
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "0"
$ws.Range("C23").NumberFormat = "General"
$ws.Range("D23").Value = 2
$ws.Range("D23").NumberFormat = "#,##0"
$ws.Range("E23").Value = -100
$ws.Range("E23").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = -75
$ws.Range("J23").Value = 22
$ws.Range("K23").Value = 4.545454545454
$ws.Range("D24").Value = 12
$ws.Range("E24").Value = 58.333333333333
$ws.Range("G24").Value = 76
$ws.Range("H24").Value = -10.526315789473
$ws.Range("I24").Value = 1082
$ws.Range("J24").Value = 1088
$ws.Range("K24").Value = -0.551470588235
$ws.Range("L24").Value = -1.096892138939
$ws.Range("M24").Value = 18.640350877193
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 16.666666666666
$ws.Range("F25").Value = 20
$ws.Range("G25").Value = 17
$ws.Range("H25").Value = 17.647058823529
$ws.Range("I25").Value = 186
$ws.Range("J25").Value = 157
$ws.Range("K25").Value = 18.471337579617
$ws.Range("L25").Value = 60.344827586206
$ws.Range("M25").Value = -21.186440677966
$ws.Range("L26").Value = 100
$ws.Range("D27").Value = 4
$ws.Range("D27").NumberFormat = "#,##0"
$ws.Range("E27").Value = -100
$ws.Range("E27").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = -60
$ws.Range("J27").Value = 44
$ws.Range("K27").Value = 4.545454545454
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "***.*"
$ws.Range("E28").NumberFormat = "General"
$ws.Range("G28").Value = 1
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "***.*"
$ws.Range("E29").NumberFormat = "General"
$ws.Range("G29").Value = 1
$ws.Range("F30").Value = 1
$ws.Range("F30").NumberFormat = "#,##0"
$ws.Range("I30").Value = 5
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 66.666666666666
